$wb = $excel.ActiveWorkbook

# --- day_frequency sheet ---
$ws = $wb.Worksheets.Item("day_frequency")
$vals = @(
    0.4997836293541847,
    0.4997847696149227,
    0.4998128481427201,
    0.4998675690634674,
    0.499948636503055,
    0.5000557545873736,
    0.5001886274423136,
    0.5003469591937658,
    0.5005304539676205,
    0.5007388158897683,
    0.5009717490860999,
    0.5012289576825055,
    0.5015101458048761,
    0.5018150175791019,
    0.5021432771310734,
    0.5024946285866814,
    0.5028687760718162,
    0.5032654237123686,
    0.5036842756342289,
    0.5041250359632878,
    0.5045874088254358,
    0.5050710983465634,
    0.5055758086525611,
    0.5061012438693197,
    0.5066471081227294,
    0.5072131055386809,
    0.5077989402430648,
    0.5084043163617715,
    0.5090289380206917,
    0.5096725093457158,
    0.5103347344627345,
    0.5110153174976381,
    0.5117139625763174,
    0.5124303738246627,
    0.5131642553685648,
    0.5139153113339141,
    0.5146832458466011,
    0.5154677630325164,
    0.5162685670175505,
    0.5170853619275941,
    0.5179178518885376,
    0.5187657410262715,
    0.5196287334666865,
    0.5205065333356729,
    0.5213988447591215,
    0.5223053718629227,
    0.523225818772967,
    0.5241598896151451,
    0.5251072885153475,
    0.5260677195994646,
    0.5270408869933871,
    0.5280264948230055,
    0.5290242472142103,
    0.530033848292892,
    0.5310550021849413,
    0.5320874130162485,
    0.5331307849127045,
    0.5341848220001995,
    0.5352492284046242,
    0.5363237082518691,
    0.5374079656678248,
    0.5385017047783818,
    0.5396046297094307,
    0.5407164445868619,
    0.541836853536566,
    0.5429655606844336,
    0.5441022701563553,
    0.5452466860782214,
    0.5463985125759226,
    0.5475574537753496,
    0.5487232138023926,
    0.5498954967829425,
    0.5510740068428895,
    0.5522584481081244,
    0.5534485247045375,
    0.5546439407580197,
    0.5558444003944611,
    0.5570496077397527,
    0.5582592669197847,
    0.5594730820604478,
    0.5606907572876324,
    0.5619119967272291,
    0.5631365045051286,
    0.5643639847472213,
    0.5655941415793978,
    0.5668266791275485,
    0.5680613015175642,
    0.5692977128753351,
    0.5705356173267521,
    0.5717747189977055,
    0.5730147220140859,
    0.5742553305017838,
    0.5754962485866899,
    0.5767371803946946,
    0.5779778300516885,
    0.5792179016835621,
    0.5804570994162059,
    0.5816951273755107,
    0.5829316896873666,
    0.5841664904776647,
    0.585399233872295,
    0.5866296239971482,
    0.5878573649781151,
    0.589082160941086,
    0.5903037160119515,
    0.5915217343166022,
    0.5927359199809286,
    0.5939459771308211,
    0.5951516098921704,
    0.596352522390867,
    0.5975484187528015,
    0.5987390031038644,
    0.5999239795699463,
    0.6011030522769375,
    0.6022759253507288,
    0.6034423029172107,
    0.6046018891022736,
    0.6057543880318083,
    0.6068995038317051,
    0.6080369406278546,
    0.6091664025461474,
    0.6102875937124741,
    0.611400218252725,
    0.6125039802927909,
    0.6135985839585621,
    0.6146837333759294,
    0.6157591326707832,
    0.6168244859690141,
    0.6178794973965126,
    0.6189238710791691,
    0.6199573111428744,
    0.620979521713519,
    0.6219902069169934,
    0.6229890708791879,
    0.6239758177259934,
    0.6249501515833003,
    0.625911776576999,
    0.6268603968329803,
    0.6277957164771346,
    0.6287174396353525,
    0.6296252704335246,
    0.6305189129975411,
    0.631398071453293,
    0.6322624499266705,
    0.6331117525435643,
    0.633945683429865,
    0.634763946711463,
    0.6355662465142489,
    0.6363522869641132,
    0.6371217721869465,
    0.6378744063086395,
    0.6386098934550823,
    0.6393279377521659,
    0.6400282433257806,
    0.6407105143018169,
    0.6413744548061655,
    0.6420197689647169,
    0.6426461609033616,
    0.6432533347479902,
    0.6438409946244932,
    0.644408844658761,
    0.6449565889766844,
    0.6454839317041539,
    0.6459905769670597,
    0.6464762288912929,
    0.6469405916027436,
    0.6473833692273026,
    0.6478042658908603,
    0.6482029857193072,
    0.648579232838534,
    0.6489327113744312,
    0.6492631254528892,
    0.6495701791997988,
    0.6498535767410502,
    0.6501130222025342,
    0.6503482197101413,
    0.650558873389762,
    0.6507446873672869,
    0.6509053657686064,
    0.6510406127196111,
    0.6511501323461917,
    0.6512336287742386,
    0.6512908061296423,
    0.6513213685382935,
    0.6513250201260826,
    0.6513014650189002,
    0.6512504073426368,
    0.651171551223183,
    0.6510646007864294,
    0.6509292601582664,
    0.6507652334645846,
    0.6505722248312744,
    0.6503499383842266,
    0.6500980782493317,
    0.64981634855248,
    0.6495044534195624,
    0.6491620969764691,
    0.6487889833490909,
    0.6483848166633182,
    0.6479493010450414,
    0.6474821406201514,
    0.6469830395145386,
    0.6464517018540934,
    0.6458878317647064,
    0.6452911333722683,
    0.6446613108026695,
    0.6439980681818006,
    0.6433011096355521,
    0.6425701392898144,
    0.6418048612704783,
    0.6410049797034343,
    0.6401701987145728,
    0.6393002224297845,
    0.6383947549749598,
    0.6374535004759894,
    0.6364761630587635,
    0.6354624468491731,
    0.6344120559731085,
    0.6333246945564602,
    0.6322000667251189,
    0.6310378766049749,
    0.629837828321919,
    0.6285996260018417,
    0.6273229737706334,
    0.6260075757541848,
    0.6246531360783862,
    0.6232593588691284,
    0.6218259482523019,
    0.6203526083537971,
    0.6188390432995047,
    0.6172849572153152,
    0.6156900542271191,
    0.6140540384608069,
    0.6123766140422692,
    0.6106574850973966,
    0.6088963557520795,
    0.6070929301322086,
    0.6052469123636743,
    0.6033580065723673,
    0.601425916884178,
    0.5994503474249969,
    0.5974310023207148,
    0.595367585697222,
    0.5932598016804091,
    0.5911073543961667,
    0.5889099479703852,
    0.5866672865289554,
    0.5843790741977676,
    0.5820450151027123,
    0.5796648133696803,
    0.577238173124562,
    0.574764798493248,
    0.5722443936016288,
    0.5696766625755949,
    0.5670613095410368,
    0.5643980386238452,
    0.5616865539499105,
    0.5589265596451234,
    0.5561177598353743,
    0.5532598586465538,
    0.5503525602045524,
    0.5473955686352606,
    0.5443885880645691,
    0.5413313226183684,
    0.538223476422549,
    0.5350647536030013,
    0.5318548582856161,
    0.5285934945962838,
    0.5252803666608948,
    0.52191517860534,
    0.5184976345555097,
    0.5150274386372945,
    0.5115042949765849,
    0.5079279076992715,
    0.5042979809312448,
    0.5006142187983954,
    0.4968763254266137,
    0.4930840049417904,
    0.489236961469816,
    0.485334899136581,
    0.481377522067976,
    0.4773645343898914,
    0.473295640228218,
    0.469170543708846,
    0.4649889489576663,
    0.4607505601005693,
    0.4564550812634454,
    0.4521022165721852,
    0.4476916701526795,
    0.4432231461308184,
    0.4386963486324928,
    0.4341109817835931,
    0.4294667497100099,
    0.4247633565376338,
    0.4200005063923551,
    0.4151779034000646,
    0.4102952516866527,
    0.40535225537801,
    0.400348618600027,
    0.3952840454785942,
    0.3901582401396023,
    0.3849709067089417,
    0.379721749312503,
    0.3744104720761767,
    0.3690367791258534,
    0.3636003745874236,
    0.3581009625867779,
    0.3525382472498068,
    0.3469119327024007,
    0.3412217230704503,
    0.3354673224798464,
    0.3296484350564791,
    0.323764764926239,
    0.3178160162150168,
    0.311801893048703,
    0.3057220995531882,
    0.2995763398543627,
    0.2933643180781174,
    0.2870857383503425,
    0.2807403047969287,
    0.2743277215437667,
    0.2678476927167464,
    0.2612999224417591,
    0.2546841148446952,
    0.2479999740514451,
    0.2412472041878994,
    0.2344255093799485,
    0.227534593753483,
    0.2205741614343935,
    0.2135439165485706,
    0.2064435632219047,
    0.1992728055802864,
    0.1920313477496062,
    0.1847188938557547,
    0.1773351480246225,
    0.1698798143821,
    0.1623525970540778,
    0.1547532001664465,
    0.1470813278450965,
    0.1393366842159186,
    0.131518973404803,
    0.1236278995376405,
    0.1156631667403216,
    0.1076244791387367,
    0.09951154085877656,
    0.09132405602633153,
    0.08306172876729218,
    0.07472426320754916,
    0.06631136347299299,
    0.05782273368951413,
    0.04925807798300319,
    0.04061710047935069,
    0.03189950530444713,
    0.02310499658418308,
    0.01423327844444916,
    0.005284055011135769,
    0.005284055011135769,
    0.005284055011135769,
    0.005284055011135769,
    0.005284055011135769,
    0.005284055011135769,
    0.005284055011135769,
    0.005284055011135769,
    0.005284055011135769,
    0.005284055011135769
)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $vals[$i]
}

# --- week_frequency sheet ---
$ws = $wb.Worksheets.Item("week_frequency")
$vals = @(
    0.4999202621011481,
    0.5010202998862483,
    0.5033098322749846,
    0.506687374447813,
    0.5110514415851893,
    0.5163005488675693,
    0.522333211475409,
    0.5290479445891642,
    0.5363432633892906,
    0.5441176830562443,
    0.5522697187704811,
    0.5606978857124567,
    0.5693006990626269,
    0.577976674001448,
    0.5866243257093754,
    0.5951421693668651,
    0.6034287201543729,
    0.6113824932523547,
    0.6189020038412665,
    0.625885767101564,
    0.632232298213703,
    0.6378401123581395,
    0.6426077247153293,
    0.6464336504657282,
    0.6492164047897921,
    0.650854502867977,
    0.6512464598807384,
    0.6502907910085324,
    0.6478860114318147,
    0.6439306363310414,
    0.6383231808866682,
    0.630962160279151,
    0.6217460896889454,
    0.6105734842965077,
    0.5973428592822935,
    0.5819527298267586,
    0.5643016111103593,
    0.5442880183135507,
    0.5218104666167891,
    0.4967674712005303,
    0.4690575472452303,
    0.4385792099313446,
    0.4052309744393295,
    0.3689113559496405,
    0.3295188696427336,
    0.2869520306990647,
    0.2411093542990891,
    0.1918893556232635,
    0.1391905498520435,
    0.08291145216588476,
    0.02424015268824354,
    0.005284055011135769
)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $vals[$i]
}

# --- month_frequency sheet ---
$ws = $wb.Worksheets.Item("month_frequency")
$vals = @(
    0.503406469767758,
    0.5237707760583156,
    0.556618728880796,
    0.5937690015890918,
    0.6268259774968072,
    0.6473040773505002,
    0.6466835087805062,
    0.618283700309273,
    0.5541537046725179,
    0.4432873880037798,
    0.2780543168069942,
    0.06455157791381295
)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $vals[$i]
}
